$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data describing LeetCode 178. Rank Scores.
# Fill cells in the same order the strings were authored (Question, then
# Link, then Notes) so the shared-string table grows in that sequence.
$ws.Range("A8").Value = "178. Rank Scores"

$ws.Range("E8").Value = "https://leetcode.com/problems/rank-scores/solutions/3861595/pandas-simple-2-step-approach-additional-knowledge-at-the-end/?envType=study-plan-v2&envId=30-days-of-pandas&lang=pythondata "
$ws.Hyperlinks.Add($ws.Range("E8"), "https://leetcode.com/problems/rank-scores/solutions/3861595/pandas-simple-2-step-approach-additional-knowledge-at-the-end/?envType=study-plan-v2&envId=30-days-of-pandas&lang=pythondata ")
$ws.Range("E8").Style = "Hyperlink"

$ws.Range("D8").Value = "Use rank method with method='dense' and descending order to rank, then drop id column and sort values by score in descending order."

# Reuse the existing "Medium" / "Data Manipulation" strings, matching the
# formatting already used on the other Medium-difficulty rows (5-7).
$ws.Range("B8").Value = "Medium"
$ws.Range("B8").Interior.Color = $ws.Range("B7").Interior.Color
$ws.Range("C8").Value = "Data Manipulation"

# Grow the worksheet table to include the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E8"))

# Match the author's final selection in the saved workbook.
[void]$ws.Range("E18").Select()
